$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 3) to the GILD "Noun" sheet, mirroring the
# existing row 2 layout: date/time serial in A, numeric metrics in B:M,
# and the "Noun" method label (shared string) in N.
$ws.Range("A3").Value = 42605.885335648149
$ws.Range("B3").Value = -28
$ws.Range("C3").Value = 53
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 81
$ws.Range("G3").Value = 21614
$ws.Range("H3").Value = 6743
$ws.Range("I3").Value = 1173
$ws.Range("J3").Value = 131
$ws.Range("K3").Value = 112
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 26
$ws.Range("N3").Value = "Noun"
